$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that is updated for every data
# row (rows 2-420) from 45190 (2023-09-21) to 45192 (2023-09-23).
for ($r = 2; $r -le 420; $r++) {
    $ws.Cells.Item($r, 3).Value = 45192
}
